$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new player row (Ja'Marr Chase) as row 6, following the existing table layout:
# A=id, B=longname, C=position, D=totalPPR, E=avgPPR, F=team, G=status, H=age,
# I=height, J=weight, K=jersey, L=notes, M..AC=week1..week17

$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "Ja'Marr Chase"
$ws.Cells.Item(6, 3).Value = "WR"
$ws.Cells.Item(6, 4).Value = 262.72000000000003
$ws.Cells.Item(6, 5).Value = 23.72
$ws.Cells.Item(6, 9).Value = "6'0"""
$ws.Cells.Item(6, 6).Value = "CIN"
$ws.Cells.Item(6, 8).Value = 24
$ws.Cells.Item(6, 10).Value = 201
$ws.Cells.Item(6, 11).Value = 1

$ws.Cells.Item(6, 13).Value = 4.0999999999999996
$ws.Cells.Item(6, 14).Value = 3.1
$ws.Cells.Item(6, 15).Value = 14.1
$ws.Cells.Item(6, 16).Value = 7.3
$ws.Cells.Item(6, 17).Value = 37.200000000000003
$ws.Cells.Item(6, 18).Value = 8
$ws.Cells.Item(6, 19).Value = 15.2
$ws.Cells.Item(6, 20).Value = 4.0999999999999996
$ws.Cells.Item(6, 21).Value = 18.399999999999999
$ws.Cells.Item(6, 22).Value = 7.2
$ws.Cells.Item(6, 23).Value = 8.1
$ws.Cells.Item(6, 24).Value = 20.6
$ws.Cells.Item(6, 25).Value = 2.9
$ws.Cells.Item(6, 26).Value = 6.4
$ws.Cells.Item(6, 27).Value = -1
$ws.Cells.Item(6, 28).Value = 4.0999999999999996
$ws.Cells.Item(6, 29).Value = 1.9

# Update the active selection on the sheet to match the author's final cursor position
$ws.Range("G12").Select()
